# The underlying data table (rows 2-177) represents a weekly price log,
# ordered from most-recent to oldest. This edit adds a new, more recent
# observation at the top of the data (row 56 in the sheet, since the
# header occupies row 1 and rows 2-55 already hold newer entries),
# pushing all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 56; this shifts old rows 56..177 down to
# 57..178 (carrying their values/formatting with them), growing the
# sheet's used range to A1:R178.
$ws.Rows("56:56").Insert()

# Populate the newly inserted row 56 with the new weekly observation.
$ws.Range("A56").Value = 4
$ws.Range("B56").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C56").Value = "Los Lagos"
$ws.Range("D56").Value = 44519
$ws.Range("E56").Value = 10
$ws.Range("F56").Value = 100112044
$ws.Range("G56").Value = "Perejil"
$ws.Range("H56").Value = "Sin especificar"
$ws.Range("I56").Value = "Primera"
$ws.Range("J56").Value = 180
$ws.Range("K56").Value = 5000
$ws.Range("L56").Value = 5000
$ws.Range("M56").Value = 5000
$ws.Range("N56").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O56").Value = "Región Metropolitana"
$ws.Range("P56").Value = 1667
$ws.Range("Q56").Value = 3
$ws.Range("R56").Value = "Hortaliza"
